# Apply "Api v3" header-label updates to the TablasTablas sheet.
# The small summary table (F3/L3/P3) and the three detail tables (row 6
# headers) are relabeled from the old Spanish field names to the new
# API v3 field names. Numeric/ID data and all other cells are untouched.
#
# NOTE on ordering: the new label strings are appended to the shared
# string table in the order they are first written, so the cells below
# are touched in the same order the target workbook introduces them
# (F3, E6, F6, G6, H6, L6, then L3) to reproduce the same shared-string
# layout as the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mini summary row (row 3) table-name label ---
$ws.Range("F3").Value = "TableTables"
# P3 ("Valores Indicadores(Sábana)") is unchanged.

# --- Header row (row 6) column labels ---
$ws.Range("E6").Value = "idindicator"
$ws.Range("F6").Value = "nametable"
$ws.Range("G6").Value = "idfield"
$ws.Range("H6").Value = "value"
$ws.Range("L6").Value = "nameindicator"
# D6/K6/O6 ("ID"), P6 ("IndicadorID"), Q6 ("Valor") and R6 ("ProcesoID") are unchanged.

# --- Mini summary row (row 3), second table-name label ---
$ws.Range("L3").Value = "Indicator"

# --- View state: zoom in a bit and move the selection/viewport ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("D3:L11").Select()
